# This workbook contains a weekly price log for "Acelga" (Swiss chard)
# at "Vega Monumental Concepción". A new weekly observation (dated
# 2022-02-11, serial 44603) is inserted as a new row right above the
# existing row for 2021-01-21 (serial 44217), which currently sits at
# row 152. All rows from the old row 152 down to the old last row (216)
# shift down by one, and the sheet's used range grows from R216 to R217.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing row at/after 152 down by one, carrying the
# formatting of row 152 (including the date-format style on column D)
# onto the freshly inserted blank row.
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with the new weekly record.
$ws.Range("A152").Value = 11
$ws.Range("B152").Value = "Vega Monumental Concepción"
$ws.Range("C152").Value = "Bíobío"
$ws.Range("D152").Value = 44603
$ws.Range("E152").Value = 8
$ws.Range("F152").Value = 100112009
$ws.Range("G152").Value = "Acelga"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 450
$ws.Range("K152").Value = 500
$ws.Range("L152").Value = 550
$ws.Range("M152").Value = 522
$ws.Range("N152").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O152").Value = "Región Metropolitana"
$ws.Range("P152").Value = 522
$ws.Range("Q152").Value = 1
$ws.Range("R152").Value = "Hortaliza"
